$d = $word.ActiveDocument

# wdReplaceOne = 1, wdFindContinue (wrap) = 1

# --- 1. Main body: "A TERE," -> "A QWER," (bold run) ---
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# --- 2-11. Header (primary header of section 1) ---
$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdrRng = $hdr.Range

# "DIRETORIA DE ENSINO REGIAO TRE" -> "...QWER"
$hdrRng.Find.Execute("TRE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# "TERE - DEP." -> "QWER - DEP."
$hdrRng.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# Address line: five "Tre" occurrences, in order -> Qwer, Qwer, Qewr, Qewr, Qwer
$hdrRng.Find.Execute("Tre", $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 1) | Out-Null
$hdrRng.Find.Execute("Tre", $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 1) | Out-Null
$hdrRng.Find.Execute("Tre", $true, $true, $false, $false, $false, $true, 1, $false, "Qewr", 1) | Out-Null
$hdrRng.Find.Execute("Tre", $true, $true, $false, $false, $false, $true, 1, $false, "Qewr", 1) | Out-Null
$hdrRng.Find.Execute("Tre", $true, $true, $false, $false, $false, $true, 1, $false, "Qwer", 1) | Out-Null

# CEP line: two "tre" occurrences -> qwer, qwer
$hdrRng.Find.Execute("tre", $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
$hdrRng.Find.Execute("tre", $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null

# Email line: "tre" -> "qwer"
$hdrRng.Find.Execute("tre", $true, $true, $false, $false, $false, $true, 1, $false, "qwer", 1) | Out-Null
